# Atualizado por script em 01-11-2023 20:45
# Adds the new match row (row 29) to the ISL 2023-2024 sheet, mirroring
# the formatting of the previous data row (row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 29
$srcRow = 28

# --- Write the values for the new row -------------------------------------
$ws.Range("A$newRow").Value = 28
$ws.Range("B$newRow").Value = "india"
$ws.Range("C$newRow").Value = "isl"
$ws.Range("D$newRow").Value = "2023-2024"
$ws.Range("E$newRow").Value = 45231.64583333334
$ws.Range("F$newRow").Value = "Jamshedpur"
$ws.Range("G$newRow").Value = 2
$ws.Range("H$newRow").Value = "Mohun Bagan"
$ws.Range("I$newRow").Value = 3
$ws.Range("J$newRow").Value = 3.62
$ws.Range("K$newRow").Value = "31/10/2023 15:30"
$ws.Range("L$newRow").Value = 4.31
$ws.Range("M$newRow").Value = "01/11/2023 15:29"
$ws.Range("N$newRow").Value = 3.45
$ws.Range("O$newRow").Value = "31/10/2023 15:30"
$ws.Range("P$newRow").Value = 3.7
$ws.Range("Q$newRow").Value = "01/11/2023 15:29"
$ws.Range("R$newRow").Value = 1.97
$ws.Range("S$newRow").Value = "31/10/2023 15:30"
$ws.Range("T$newRow").Value = 1.83
$ws.Range("U$newRow").Value = "01/11/2023 15:29"
$ws.Range("V$newRow").Value = "https://www.betexplorer.com/football/india/isl/jamshedpur-mohun-bagan/EXDAFS8d/"

# --- Mirror formatting (styles) from the row above -------------------------
$ws.Range("A$srcRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E$srcRow").Copy()
$ws.Range("E$newRow").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
